$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rich text header edits ---
$a8 = $ws.Range("A8")
$r1 = $a8.Characters(21,2)
$r1.Text = "36"
$r1.Font.Name = "Andale WT"
$r1.Font.Size = 10

$c9 = $ws.Range("C9")
$d1 = $c9.Characters(27,9)
$d1.Text = "9/2/2024"
$d1.Font.Name = "Andale WT"
$d1.Font.Size = 10
$d2 = $c9.Characters(46,8)
$d2.Text = "9/8/2024"
$d2.Font.Name = "Andale WT"
$d2.Font.Size = 10

# --- Numeric value updates ---
$ws.Range("N15").Value = -56.666666666666
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 128
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3.225806451612
$ws.Range("M16").Value = -47.107438016528
$ws.Range("N16").Value = -85.487528344671
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 215
$ws.Range("J17").Value = 213
$ws.Range("K17").Value = 0.93896713615
$ws.Range("L17").Value = 8.040201005025
$ws.Range("M17").Value = 50.34965034965
$ws.Range("N17").Value = -46.650124069478
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -24.137931034482
$ws.Range("I18").Value = 230
$ws.Range("J18").Value = 187
$ws.Range("K18").Value = 22.994652406417
$ws.Range("L18").Value = 11.650485436893
$ws.Range("M18").Value = -25.324675324675
$ws.Range("N18").Value = -75.453575240128
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 7.54716981132
$ws.Range("I19").Value = 496
$ws.Range("J19").Value = 523
$ws.Range("K19").Value = -5.162523900573
$ws.Range("L19").Value = 8.061002178649
$ws.Range("M19").Value = 52.615384615384
$ws.Range("N19").Value = 35.890410958904
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 81
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = -31.355932203389
$ws.Range("L20").Value = -34.677419354838
$ws.Range("M20").Value = -22.115384615384
$ws.Range("N20").Value = -86.43216080402
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 57.692307692307
$ws.Range("F21").Value = 129
$ws.Range("G21").Value = 132
$ws.Range("H21").Value = -2.272727272727
$ws.Range("I21").Value = 1166
$ws.Range("J21").Value = 1181
$ws.Range("K21").Value = -1.270110076206
$ws.Range("L21").Value = 3.552397868561
$ws.Range("M21").Value = 3.460514640638
$ws.Range("N21").Value = -63.900928792569
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 117
$ws.Range("J23").Value = 141
$ws.Range("K23").Value = -17.021276595744
$ws.Range("L23").Value = 6.363636363636
$ws.Range("M23").Value = 20.61855670103
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 29.411764705882
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 27.397260273972
$ws.Range("I24").Value = 777
$ws.Range("J24").Value = 731
$ws.Range("K24").Value = 6.292749658002
$ws.Range("L24").Value = -9.334889148191
$ws.Range("M24").Value = -7.829181494661
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -26.086956521739
$ws.Range("I25").Value = 244
$ws.Range("J25").Value = 105
$ws.Range("K25").Value = 132.380952380952
$ws.Range("L25").Value = 25.773195876288
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -23.076923076923
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 361
$ws.Range("J26").Value = 349
$ws.Range("K26").Value = 3.438395415472
$ws.Range("L26").Value = 6.489675516224
$ws.Range("M26").Value = 9.063444108761
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = -15.78947368421
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 60
$ws.Range("L28").Value = 36.363636363636
$ws.Range("L29").Value = -11.111111111111
$ws.Range("L30").Value = -22.222222222222
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = 72.727272727272
$ws.Range("L31").Value = 26.666666666666
$ws.Range("F33").Value = 2
$ws.Range("I33").Value = 5
$ws.Range("K33").Value = 150
$ws.Range("L33").Value = 150

# --- Numeric -> Text "0" conversions (style 14) ---
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)

# --- Numeric -> Text "***.*" conversions (style 14) ---
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

# --- Text -> Numeric conversions ---
$ws.Range("D31").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E31").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
